$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D text values that look numeric (prices with single decimal
# point, or many leading zeros) are written back as literal text, matching
# the inlineStr cells already used throughout the sheet, instead of being
# auto-coerced to floating point numbers by Excel (which would silently drop
# meaningful trailing zeros, e.g. "0.00001090" -> 0.0000109).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.558.86"
$ws.Range("E2").Value = "  +1.80%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.826.63"
$ws.Range("E3").Value = "  +1.84%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.34"
$ws.Range("E5").Value = "  +0.05%  "

$ws.Range("E6").Value = "  +0.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5405"
$ws.Range("E7").Value = "  +0.77%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4004"
$ws.Range("E8").Value = "  +6.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07753"
$ws.Range("E9").Value = "  +4.33%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.121"
$ws.Range("E10").Value = "  +2.59%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.99"
$ws.Range("E11").Value = "  +0.38%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.29"
$ws.Range("E12").Value = "  +3.54%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.346"
$ws.Range("E13").Value = "  +3.63%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.627"
$ws.Range("E14").Value = "  +5.28%  "

$ws.Range("E15").Value = "  -0.01%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.828.82"
$ws.Range("E16").Value = "  +2.09%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001090"
$ws.Range("E17").Value = "  +2.97%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "89.92"
$ws.Range("E18").Value = "  +1.02%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06588"
$ws.Range("E19").Value = "  +1.32%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.75"
$ws.Range("E20").Value = "  +2.86%  "

$ws.Range("E21").Value = "  +0.10%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.077"
$ws.Range("E22").Value = "  +3.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.575.24"
$ws.Range("E23").Value = "  +1.85%  "

$ws.Range("E24").Value = "  +0.46%  "

$ws.Range("E25").Value = "  +8.35%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.16"
$ws.Range("E26").Value = "  +1.72%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.84"
$ws.Range("E27").Value = "  +2.67%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.454"
$ws.Range("E28").Value = "  +6.63%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.039.85"
$ws.Range("E29").Value = "  +2.31%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.45"
$ws.Range("E30").Value = "  +2.64%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.135"
$ws.Range("E31").Value = "  +1.64%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1121"
$ws.Range("E32").Value = "  +5.78%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.701"
$ws.Range("E33").Value = "  +2.58%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07415"
$ws.Range("E34").Value = "  +14.22%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.647"
$ws.Range("E35").Value = "  -0.43%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2259"
$ws.Range("E36").Value = "  +0.46%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02360"
$ws.Range("E37").Value = "  +3.09%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.964"
$ws.Range("E38").Value = "  +5.95%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.218"
$ws.Range("E39").Value = "  +4.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.42"
$ws.Range("E40").Value = "  +2.56%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6308"
$ws.Range("E41").Value = "  +1.84%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.191"
$ws.Range("E42").Value = "  +1.34%  "

$ws.Range("E43").Value = "  +0.03%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.405"
$ws.Range("E44").Value = "  -3.15%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.56"
$ws.Range("E45").Value = "  +2.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5906"
$ws.Range("E46").Value = "  +2.16%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.712"
$ws.Range("E47").Value = "  +1.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.38"
$ws.Range("E48").Value = "  +0.23%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.003"
$ws.Range("E49").Value = "  +3.95%  "

$ws.Range("E50").Value = "  +0.88%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06919"
$ws.Range("E51").Value = "  +1.47%  "
